$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.820.77"
$ws.Range("E2").Value = "  +2.56%  "
$ws.Range("D3").Value = "3.807.99"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "698.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.14%  "
$ws.Range("D7").Value = "3.806.62"
$ws.Range("E7").Value = "  +1.24%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.528"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.162"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.27"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.459"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000256"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +7.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.31"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.03%  "
$ws.Range("D15").Value = "4.455.01"
$ws.Range("E15").Value = "  +1.47%  "
$ws.Range("D16").Value = "3.791.47"
$ws.Range("E16").Value = "  +0.84%  "
$ws.Range("D17").Value = "70.832.31"
$ws.Range("E17").Value = "  +2.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.82"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.63%  "
$ws.Range("E20").Value = "  +0.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +16.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "480.23"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.712"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.47%  "
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.16"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.44"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.86%  "
$ws.Range("D29").Value = "3.959.86"
$ws.Range("E29").Value = "  +1.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E31").Value = "  +14.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.29"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.51"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.34%  "
$ws.Range("E34").Value = "  +11.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "29.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.22"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.05%  "
$ws.Range("E37").Value = "  +0.34%  "
$ws.Range("E38").Value = "  +2.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.42"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +12.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.000330"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +24.30%  "
$ws.Range("E43").Value = "  +1.94%  "
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "162.63"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "44.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "48.73"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.54%  "
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.39"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.45%  "
$ws.Range("B50").Value = "TheGraph"
$ws.Range("C50").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.300"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.56"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.41%  "
